# "Busqueda aleatoria y clic en articulo" -- the inventory search picked a
# different product ("televisor" instead of "taladro") and the user then
# clicked on that item, moving the active selection to it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the product name in the inventory list.
$ws.Range("A2").Value = "televisor"

# Click on the article that was just found/updated.
$ws.Range("A2").Select()
